# Revert "pre alfa, se puede asignar prioridad a los ramos y a las secciones,
# se muestra todas las opciones y la de mayor peso, limitado a 6 ramos por semestre"
#
# Net effect on sheet data: the per-section "Electivo Profesional-N" labels
# collapse back to the single generic "Electivo Profesional" label, and two
# unrelated data corrections (D27, F55) revert to their prior values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the numbered elective labels back to the generic one.
$ws.Range("C45").Value2 = "Electivo Profesional"
$ws.Range("C46").Value2 = "Electivo Profesional"
$ws.Range("C47").Value2 = "Electivo Profesional"
$ws.Range("C49").Value2 = "Electivo Profesional"
$ws.Range("C50").Value2 = "Electivo Profesional"
$ws.Range("C51").Value2 = "Electivo Profesional"
$ws.Range("C52").Value2 = "Electivo Profesional"
$ws.Range("C54").Value2 = "Electivo Profesional"

# Other two data reverts.
$ws.Range("D27").Value2 = "31, 35"
$ws.Range("F55").Value2 = "13"

# Restore the prior view state (scroll position / active selection).
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("J44").Select()
